# ContosoLearn Competitor SWOT - localize document body text to Indonesian.
# Each paragraph's existing runs (bold "label:" run + normal body run) keep
# their original character formatting; only the underlying text is swapped,
# via Find/Replace scoped to that paragraph's Range so nothing outside the
# matched text is touched.

$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range
    $rng.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)
}

# Paragraph 1: title
Replace-InParagraph 1 "ContosoLearn Competitor SWOT" "SWOT Kompetitor ContosoLearn"

# Paragraph 2: "Fabrikam Learning:" heading
Replace-InParagraph 2 "Fabrikam Learning:" "Pembelajaran Fabrikam:"

# Paragraph 3: Strengths (Fabrikam)
Replace-InParagraph 3 "Strengths:" "Kekuatan:"
Replace-InParagraph 3 " Fabrikam Learning provides a comprehensive set of analytics and reporting tools. It ensures the continuous monitoring of teaching and learning activities, as well as pinpointing problematic areas that need to be addressed." " Pembelajaran Fabrikam menyediakan serangkaian alat analitik dan pelaporan yang komprehensif. Ini memastikan pemantauan berkelanjutan terhadap kegiatan belajar mengajar, serta menentukan area bermasalah yang perlu ditangani."

# Paragraph 4: Weaknesses (Fabrikam)
Replace-InParagraph 4 "Weaknesses:" "Kelemahan:"
Replace-InParagraph 4 " While Fabrikam Learning has robust reporting capabilities, it might be overwhelming for some users due to its comprehensive nature." " Meskipun Pembelajaran Fabrikam memiliki kemampuan pelaporan yang kuat, beberapa pengguna mungkin akan kewalahan karena sifatnya yang komprehensif."

# Paragraph 5: Opportunities (Fabrikam)
Replace-InParagraph 5 "Opportunities:" "Peluang:"
Replace-InParagraph 5 " There is a growing demand for personalized learning experiences and data-driven recommendations. Fabrikam Learning can leverage its robust analytics and reporting tools to meet this demand." " Ada peningkatan permintaan untuk pengalaman pembelajaran yang dipersonalisasi dan rekomendasi berbasis data. Pembelajaran Fabrikam dapat memanfaatkan analitik dan alat pelaporan yang kuat untuk memenuhi permintaan ini."

# Paragraph 6: Threats (Fabrikam)
Replace-InParagraph 6 "Threats:" "Ancaman:"
Replace-InParagraph 6 " The eLearning market is highly competitive with many players offering similar features. Fabrikam Learning needs to continuously innovate to stay ahead." " Pasar eLearning sangat kompetitif dengan banyaknya pemain yang menawarkan fitur serupa. Pembelajaran Fabrikam perlu terus berinovasi untuk tetap unggul."

# Paragraph 7: "AdatumLearn:" heading (text unchanged)

# Paragraph 8: Strengths (AdatumLearn)
Replace-InParagraph 8 "Strengths:" "Kekuatan:"
Replace-InParagraph 8 " AdatumLearn offers courses on business analysis techniques such as MOST and SWOT. This shows their commitment to providing valuable content to their users." " AdatumLearn menawarkan kursus tentang teknik analisis bisnis seperti MOST dan SWOT. Ini menunjukkan komitmen mereka untuk memberikan konten yang berharga bagi penggunanya."

# Paragraph 9: Weaknesses (AdatumLearn)
Replace-InParagraph 9 "Weaknesses:" "Kelemahan:"
Replace-InParagraph 9 " The information provided in their courses is a compilation of third-party generated information. This might not be as valuable as original content." " Informasi yang diberikan dalam kursus mereka merupakan kompilasi informasi yang dihasilkan pihak ketiga. Ini mungkin tidak seberharga konten asli."

# Paragraph 10: Opportunities (AdatumLearn)
Replace-InParagraph 10 "Opportunities:" "Peluang:"
Replace-InParagraph 10 " AdatumLearn can create more original content to provide unique value to their users. They can also expand their course offerings to cover more topics." " AdatumLearn dapat membuat lebih banyak konten asli untuk memberikan nilai unik bagi penggunanya. Mereka juga dapat memperluas penawaran kursusnya untuk mencakup lebih banyak topik."

# Paragraph 11: Threats (AdatumLearn)
# NB: the trailing `"` is left untouched by the Find/Replace (it is not
# included in either the search or replacement text) so Word's smart-quote
# autocorrect never sees/rewrites it into a curly quote.
Replace-InParagraph 11 "Threats:" "Ancaman:"
Replace-InParagraph 11 " Like Fabrikam Learning, AdatumLearn also faces stiff competition in the eLearning market. They need to continuously improve their offerings to stay competitive." " Seperti Pembelajaran Fabrikam, AdatumLearn juga menghadapi persaingan ketat di pasar eLearning. Mereka perlu meningkatkan penawarannya secara terus-menerus agar tetap kompetitif."
